$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G7").Value = 1.8
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 2.5
$ws.Range("L7").Value = 5
$ws.Range("AD7").Value = 7
$ws.Range("AK7").Value = 51
$ws.Range("AX7").Value = 26
$ws.Range("AZ7").Value = 101
$ws.Range("H11").Value = 3.85
$ws.Range("J11").Value = 1.91
$ws.Range("K11").Value = 2.18
$ws.Range("L11").Value = 7.8
$ws.Range("N11").Value = 6.1
$ws.Range("P11").Value = 2.65
$ws.Range("U11").Value = 2.52
$ws.Range("V11").Value = 1.47
$ws.Range("X11").Value = 5.3
$ws.Range("AB11").Value = 45
$ws.Range("AC11").Value = 6.1
$ws.Range("AD11").Value = 8
$ws.Range("AE11").Value = 29
$ws.Range("AH11").Value = 15
$ws.Range("AO11").Value = 6.4
$ws.Range("AQ11").Value = 19.5
$ws.Range("AS11").Value = 350
$ws.Range("AW11").Value = 9
$ws.Range("AX11").Value = 55
$ws.Range("G13").Value = 4.5
$ws.Range("I13").Value = 1.73
$ws.Range("Q13").Value = 1.65
$ws.Range("R13").Value = 2.2
$ws.Range("S13").Value = 1.33
$ws.Range("T13").Value = 3.25
$ws.Range("U13").Value = 1.67
$ws.Range("V13").Value = 2.1
$ws.Range("AA13").Value = 34
$ws.Range("AE13").Value = 15
$ws.Range("AH13").Value = 8.5
$ws.Range("AI13").Value = 9
$ws.Range("AN13").Value = 6.5
$ws.Range("AO13").Value = 23
$ws.Range("AT13").Value = 3.25
$ws.Range("AV13").Value = 51
$ws.Range("AZ13").Value = 26
$ws.Range("G16").Value = 1.7
$ws.Range("I16").Value = 4.5
$ws.Range("J16").Value = 2.5
$ws.Range("L16").Value = 5.5
$ws.Range("N16").Value = 7.5
$ws.Range("U16").Value = 2.2
$ws.Range("V16").Value = 1.62
$ws.Range("AC16").Value = 7.5
$ws.Range("AL16").Value = 41
$ws.Range("AU16").Value = 9.5
$ws.Range("AX16").Value = 29
$ws.Range("G18").Value = 2.25
$ws.Range("I18").Value = 3.1
$ws.Range("J18").Value = 2.88
$ws.Range("L18").Value = 3.5
$ws.Range("N18").Value = 13
$ws.Range("Q18").Value = 1.75
$ws.Range("R18").Value = 2.05
$ws.Range("W18").Value = 9.5
$ws.Range("X18").Value = 12
$ws.Range("AJ18").Value = 11
$ws.Range("AW18").Value = 5
$ws.Range("G20").Value = 1.65
$ws.Range("H20").Value = 4
$ws.Range("I20").Value = 4.2
$ws.Range("J20").Value = 2.2
$ws.Range("L20").Value = 4.5
$ws.Range("X20").Value = 9
$ws.Range("Z20").Value = 13
$ws.Range("AD20").Value = 8
$ws.Range("AI20").Value = 23
$ws.Range("AO20").Value = 8.5
$ws.Range("BA20").Value = 81
$ws.Range("G23").Value = 4.4
$ws.Range("H23").Value = 3.75
$ws.Range("I23").Value = 1.75
$ws.Range("J23").Value = 4.6
$ws.Range("P23").Value = 3.45
$ws.Range("Q23").Value = 1.85
$ws.Range("R23").Value = 1.93
$ws.Range("U23").Value = 1.78
$ws.Range("V23").Value = 1.93
$ws.Range("W23").Value = 12
$ws.Range("X23").Value = 27
$ws.Range("Y23").Value = 15.5
$ws.Range("Z23").Value = 80
$ws.Range("AD23").Value = 7.5
$ws.Range("AE23").Value = 16.5
$ws.Range("AF23").Value = 80
$ws.Range("AH23").Value = 7
$ws.Range("AI23").Value = 8.75
$ws.Range("AK23").Value = 14.5
$ws.Range("AM23").Value = 28
$ws.Range("AN23").Value = 6.1
$ws.Range("AT23").Value = 2.8
$ws.Range("AU23").Value = 7.5
$ws.Range("AW23").Value = 3.6
$ws.Range("AX23").Value = 8.5
$ws.Range("AY23").Value = 18
$ws.Range("AZ23").Value = 29
